$d = $word.ActiveDocument

# Locate a literal piece of text in the document and return the Range that
# spans exactly that text (collapsed to the match after Find.Execute).
function Find-TextRange($doc, $searchText) {
    $rng = $doc.Content
    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $rng
}

# Replace the text occupying [rangeStart, rangeStart + oldLength) with
# $newText, laid out as a sequence of separate runs - one per entry in
# $segments (which must concatenate to $newText) - all sharing the run's
# original formatting.
function Split-RangeIntoRuns($doc, $rangeStart, $oldLength, $newText, $segments) {
    # Replace the whole old text with the new text first; this lands as a
    # single run using the formatting that was already present there.
    $rangeEnd = $rangeStart + $oldLength
    $full = $doc.Range($rangeStart, $rangeEnd)
    $full.Text = $newText

    # Compute character boundaries for each segment within the new text.
    $bounds = @(0)
    $pos = 0
    foreach ($seg in $segments) {
        $pos = $pos + $seg.Length
        $bounds += $pos
    }

    # Force Word to keep each segment as its own run by toggling a formatting
    # property on and back off over each sub-range. This "touches" the run so
    # it does not get coalesced back together with its neighbours even though
    # every segment ends up with identical formatting in the end.
    for ($i = 0; $i -lt $bounds.Length - 1; $i++) {
        $a = $rangeStart + $bounds[$i]
        $b = $rangeStart + $bounds[$i + 1]
        $seg = $doc.Range($a, $b)
        $seg.Bold = 1
        $seg.Bold = 0
    }
}

# --- Start Date: "12 - Aug - 2020" -> "14 - 02 - 2021" ---
$oldText1 = "12 - Aug - 2020"
$found1 = Find-TextRange $d $oldText1
$startPos1 = $found1.Start
$oldLen1 = $oldText1.Length
$segments1 = @("1", "4", " - ", "02", " - 202", "1")
$newText1 = [string]::Join("", $segments1)
Split-RangeIntoRuns $d $startPos1 $oldLen1 $newText1 $segments1

# --- End Date: "15 - Dec - 2020" -> "02 - 05- 2021" ---
$oldText2 = "15 - Dec - 2020"
$found2 = Find-TextRange $d $oldText2
$startPos2 = $found2.Start
$oldLen2 = $oldText2.Length
$segments2 = @("02", " - ", "05", "- 202", "1")
$newText2 = [string]::Join("", $segments2)
Split-RangeIntoRuns $d $startPos2 $oldLen2 $newText2 $segments2
